# Apply cryptos.xlsx price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.846.89'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '1.630.17'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.67'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5065'
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2573'
$ws.Range('E8').Value = '  +0.49%  '
$ws.Range('E9').Value = '  -0.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.46'
$ws.Range('E10').Value = '  -0.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07763'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.248'
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('D13').Value = '1.637.72'
$ws.Range('E13').Value = '  -0.36%  '
$ws.Range('D14').Value = '1.854.09'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5500'
$ws.Range('E15').Value = '  +1.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.70'
$ws.Range('E16').Value = '  -0.67%  '
$ws.Range('D17').Value = '0.0₅7638'
$ws.Range('E17').Value = '  -2.18%  '
$ws.Range('D18').Value = '25.874.86'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.003'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.75'
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.414'
$ws.Range('E21').Value = '  -0.43%  '
$ws.Range('E22').Value = '  -0.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.029'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('E25').Value = '  +2.02%  '
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('E27').Value = '  +3.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.772'
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.58'
$ws.Range('E29').Value = '  -0.78%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.238'
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04886'
$ws.Range('E31').Value = '  -1.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.242'
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.191'
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('E34').Value = '  +0.50%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.370'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8956'
$ws.Range('E36').Value = '  +0.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5529'
$ws.Range('E37').Value = '  +2.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.540'
$ws.Range('E38').Value = '  -1.46%  '
$ws.Range('D39').Value = '1.123.77'
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01550'
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.001'
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.581'
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7940'
$ws.Range('E43').Value = '  -2.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.19'
$ws.Range('E44').Value = '  -2.19%  '
$ws.Range('D45').Value = '0.0₈118'
$ws.Range('E45').Value = '  -4.54%  '
$ws.Range('D46').Value = '1.763.44'
$ws.Range('E46').Value = '  -0.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4442'
$ws.Range('E47').Value = '  -2.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.005'
$ws.Range('E48').Value = '  +0.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.70'
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('E50').Value = '  +1.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.587'
$ws.Range('E51').Value = '  +3.18%  '
